# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6867
$wsExhibit.Range("F3").Value = 51
$wsExhibit.Range("F5").Value = 31
$wsExhibit.Range("F6").Value = 1068
$wsExhibit.Range("F7").Value = 154
$wsExhibit.Range("F8").Value = 5

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6867
$wsAll.Range("F3").Value = 51
$wsAll.Range("F5").Value = 31
$wsAll.Range("F6").Value = 1068
$wsAll.Range("F7").Value = 154
$wsAll.Range("F9").Value = 5
